$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor values change
$ws.Range("B3").Value = 0.009183249939904782
$ws.Range("C3").Value = 0.01005644336237487
$ws.Range("D3").Value = 0.01096883779990067

# Row 4 - label changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.008486287632614541
$ws.Range("C4").Value = 0.008424325166877212
$ws.Range("D4").Value = 0.008501509040853873

# Row 5 - label changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.009254582960705032
$ws.Range("C5").Value = 0.009179242514751428
$ws.Range("D5").Value = 0.008666441091175241
